# Clean up the recruitment-number codes in column A and move the
# selection to A5, matching the cleaned-up worksheet state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shorten the KRS-style reference numbers in A3 and A4.
$ws.Range("A3").Value = "0000142"
$ws.Range("A4").Value = "00004"

# Leave the active selection on A5.
[void]$ws.Range("A5").Select()
